$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Unprotect()

$ws.Range("A59").Value = "***CONFIDENTIAL***: For one-on-one client use only. Not approved for distribution.`nModel holdings provided as of 2021-04-23 for illustrative purposes only and are subject to change."

$ws.Range("D2").Value = 0.02070148208652058
$ws.Range("E2").Value = 0.01125511596180084
$ws.Range("D3").Value = 0.01853692446823597
$ws.Range("E3").Value = 0.009794319294809117
$ws.Range("D4").Value = 0.01964684021994165
$ws.Range("E4").Value = -0.002957121734844748
$ws.Range("D5").Value = 0.01983021234249717
$ws.Range("E5").Value = 0.05171870072532347
$ws.Range("D6").Value = 0.01946770441704913
$ws.Range("E6").Value = 0.02113902014424274
$ws.Range("D7").Value = 0.01942796370402005
$ws.Range("E7").Value = 0.01402805611222435
$ws.Range("D8").Value = 0.01996597630122088
$ws.Range("E8").Value = 0.0107705053852527
$ws.Range("D9").Value = 0.02006966526821048
$ws.Range("E9").Value = 0.02500804117079447
$ws.Range("D10").Value = 0.01886735740194988
$ws.Range("E10").Value = 0.01847575057736717
$ws.Range("D11").Value = 0.01949614827764355
$ws.Range("E11").Value = 0.0442340524600342
$ws.Range("D12").Value = 0.01910096017193809
$ws.Range("E12").Value = 0.01498637602179853
$ws.Range("D13").Value = 0.02108396123323703
$ws.Range("E13").Value = 0.004879635653871395
$ws.Range("D14").Value = 0.02014793631722207
$ws.Range("E14").Value = -0.003223997757218933
$ws.Range("D15").Value = 0.0191368680243197
$ws.Range("E15").Value = -0.002319109461966562
$ws.Range("D16").Value = 0.01769127437169831
$ws.Range("E16").Value = 0.01808479098725191
$ws.Range("D17").Value = 0.0178187674206038
$ws.Range("E17").Value = 0.04981320049813198
$ws.Range("D18").Value = 0.01647121430875468
$ws.Range("E18").Value = -0.01521126760563374
$ws.Range("D19").Value = 0.01470547592731524
$ws.Range("E19").Value = 0.02203108495548545
$ws.Range("D20").Value = 0.02244442503358292
$ws.Range("E20").Value = -0.003954700701060854
$ws.Range("D21").Value = 0.0216197548058525
$ws.Range("E21").Value = 0.01461202552905605
$ws.Range("D22").Value = 0.02114972505276739
$ws.Range("E22").Value = 0.009337860780984863
$ws.Range("D23").Value = 0.020376697579024
$ws.Range("E23").Value = 0.02361152361152374
$ws.Range("D24").Value = 0.01873260209076496
$ws.Range("E24").Value = 0.01173810036614253
$ws.Range("D25").Value = 0.01893998002474414
$ws.Range("E25").Value = 0.01584867075664631
$ws.Range("D26").Value = 0.01996516938318984
$ws.Range("E26").Value = 0.01333737496210996
$ws.Range("D27").Value = 0.01875398541858772
$ws.Range("E27").Value = 0.01335972290945064
$ws.Range("D28").Value = 0.01995125004715427
$ws.Range("E28").Value = 0.01801801801801806
$ws.Range("D29").Value = 0.01698158996339215
$ws.Range("E29").Value = 0.007840342124020072
$ws.Range("D30").Value = 0.01303898846369464
$ws.Range("E30").Value = 0.01327433628318597
$ws.Range("D31").Value = 0.009694111495495481
$ws.Range("E31").Value = 0.02257829570284042
$ws.Range("D32").Value = 0.01766666337175137
$ws.Range("E32").Value = -0.04471544715447151
$ws.Range("D33").Value = 0.01988467930959286
$ws.Range("E33").Value = -0.0305972344807296
$ws.Range("D34").Value = 0.02064701511942488
$ws.Range("E34").Value = -0.01033707865168543
$ws.Range("D35").Value = 0.01794343625640063
$ws.Range("E35").Value = -0.02662229617304501
$ws.Range("D36").Value = 0.02016165392374988
$ws.Range("E36").Value = 0.003722084367245637
$ws.Range("D37").Value = 0.01810724061670324
$ws.Range("E37").Value = -0.01590909090909087
$ws.Range("D38").Value = 0.02106963843818594
$ws.Range("E38").Value = -0.006319115323854319
$ws.Range("D39").Value = 0.02273895011491521
$ws.Range("E39").Value = 0.02015613910574876
$ws.Range("D40").Value = 0.01931761766327432
$ws.Range("E40").Value = 0.01984126984126977
$ws.Range("D41").Value = 0.02122839956079451
$ws.Range("E41").Value = 0.003192945111753165
$ws.Range("D42").Value = 0.0197200680312592
$ws.Range("E42").Value = 0.009329446064139768
$ws.Range("D43").Value = 0.01995165350616979
$ws.Range("E43").Value = 0.0122038765254846
$ws.Range("D44").Value = 0.01942675332697347
$ws.Range("E44").Value = 0.01886792452830188
$ws.Range("D45").Value = 0.01882136307418018
$ws.Range("E45").Value = 0.02947481243301175
$ws.Range("D46").Value = 0.01938539877788229
$ws.Range("E46").Value = 0.01808608058608052
$ws.Range("D47").Value = 0.01984715762114916
$ws.Range("E47").Value = -0.005122732123799434
$ws.Range("D48").Value = 0.01856718389440025
$ws.Range("E48").Value = 0.01108213820078241
$ws.Range("D49").Value = 0.01691199328321431
$ws.Range("E49").Value = 0.04526748971193406
$ws.Range("D50").Value = 0.01798136140385985
$ws.Range("E50").Value = 0.006865912762520132
$ws.Range("D51").Value = 0.01732755606920371
$ws.Range("E51").Value = 0.01121718377088299
$ws.Range("D52").Value = 0.01775623127319763
$ws.Range("E52").Value = 0.01295160190865707
$ws.Range("D53").Value = 0.01663945671822806
$ws.Range("E53").Value = 0.02514427040395706
$ws.Range("D54").Value = 0.007495259860891366
$ws.Range("E54").Value = 0.006055712555510739
$ws.Range("D55").Value = 0.007588257163969577
$ws.Range("E55").Value = -0.009570395576350532
$ws.Range("D56").Value = 0.9999999999999999
$ws.Range("E56").Value = 0.01036375259651123

$ws.Protect()
